# "Made multi-matrix hardware fixed software support"
#
# The "Font Workshop" sheet is a pixel-grid glyph editor: C4:J11 hold the
# 0/1 pixel bits for the glyph currently being edited, L4:L11 packs each
# row of 8 bits into a byte via a formula, and L13 concatenates the byte
# values into a "{b0,b1,...,b7}" string that gets copied into the
# matching AA/AD column cell for the letter being defined.
#
# This edit redraws the glyph in the C4:J11 grid (new pixel pattern),
# re-enters/fills the L4:L11 byte formula (which Excel re-serialises as
# a shared formula group), commits the resulting byte string into AD14
# (the "J" glyph slot), and also hand-edits AA16 (the "L" glyph slot)
# with a corrected byte string. The active selection is moved to L16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Font Workshop")

# ---------------------------------------------------------------------
# 1. Redraw the glyph pixel grid (C4:J11). Values are the 0/1 "on/off"
#    pixel bits; columns C..J are bits 0..7 (C=LSB, J=MSB) of the packed
#    byte computed in column L.
# ---------------------------------------------------------------------
$pixels = @{
    4  = @(0,0,0,0,0,0,0,0)
    5  = @(0,0,0,0,0,0,0,0)
    6  = @(0,0,0,1,0,0,0,0)
    7  = @(0,0,0,0,0,0,0,0)
    8  = @(0,0,0,1,0,0,0,0)
    9  = @(0,0,0,1,0,0,0,0)
    10 = @(0,1,0,1,0,0,0,0)
    11 = @(0,0,1,0,0,0,0,0)
}

$cols = @("C","D","E","F","G","H","I","J")
foreach ($row in 4..11) {
    $vals = $pixels[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 2. Re-enter the byte-packing formula across L4:L11 (this is what
#    turns the eight previously-independent <f> formulas into a single
#    shared formula group, matching a fill-down of L4 through L11).
# ---------------------------------------------------------------------
$ws.Range("L4:L11").Formula = "=C4+2*D4+4*E4+8*F4+G4*16+H4*32+I4*64+J4*128"

# ---------------------------------------------------------------------
# 3. Hand-correct AA16 (the "L" slot) with its new byte-string value.
# ---------------------------------------------------------------------
$ws.Range("AA16").Value = "{0,4,4,4,4,4,4,28}"

# ---------------------------------------------------------------------
# 4. Commit the freshly computed glyph byte-string into AD14 (the "J"
#    slot in the lowercase/alt glyph table) -- this mirrors L13, the
#    CONCATENATE of L4:L11 into "{...}" form.
# ---------------------------------------------------------------------
$ws.Range("AD14").Value = "{0,0,8,0,8,8,10,4}"

# ---------------------------------------------------------------------
# 5. Move the active selection to L16.
# ---------------------------------------------------------------------
$ws.Range("L16").Select() | Out-Null
